{"js": "// Replace each arithmetic expression in the worksheet table with its\n// updated value. Every \"old\" expression text is unique in the document,\n// so we can safely locate each one with a literal search and replace\n// just that matched range -- this preserves the surrounding paragraph\n// and run formatting (font, size, alignment, etc.) instead of rewriting\n// whole paragraphs/cells.\nconst replacements = [\n  [\"2+39=\", \"66+2=\"],\n  [\"81-14=\", \"36+0=\"],\n  [\"61+29=\", \"42-0=\"],\n  [\"26+31=\", \"94-35=\"],\n  [\"15-7=\", \"48-24=\"],\n  [\"62-5=\", \"15-1=\"],\n  [\"70-14=\", \"19+37=\"],\n  [\"62+12=\", \"43-1=\"],\n  [\"84+10=\", \"44+9=\"],\n  [\"31-18=\", \"45+34=\"],\n  [\"25-23=\", \"93-74=\"],\n  [\"6+83=\", \"26+62=\"],\n  [\"88-69=\", \"68+14=\"],\n  [\"53+11=\", \"24+7=\"],\n  [\"69+12=\", \"46-15=\"],\n  [\"7+86=\", \"3+8=\"],\n  [\"59+4=\", \"38+41=\"],\n  [\"40+39=\", \"32+20=\"],\n  [\"76-28=\", \"82-76=\"],\n  [\"27+27=\", \"93-26=\"],\n  [\"5+78=\", \"49-11=\"],\n  [\"10+55=\", \"73+26=\"],\n  [\"60+29=\", \"81-40=\"],\n  [\"73-69=\", \"41+20=\"],\n  [\"41-30=\", \"72+8=\"],\n  [\"40+53=\", \"36-24=\"],\n  [\"5+3=\", \"63+16=\"],\n  [\"67-48=\", \"1+55=\"],\n  [\"4+81=\", \"92-86=\"],\n  [\"93-21=\", \"74-4=\"],\n  [\"77+6=\", \"6+3=\"],\n  [\"31-23=\", \"2+63=\"],\n  [\"18+26=\", \"72+24=\"],\n  [\"70-54=\", \"86-51=\"],\n  [\"18+60=\", \"16+43=\"],\n  [\"82-62=\", \"80-63=\"],\n  [\"56-46=\", \"76+8=\"],\n  [\"90-66=\", \"57+8=\"],\n  [\"58+30=\", \"74-25=\"],\n  [\"67+17=\", \"85-59=\"],\n  [\"71+22=\", \"3+86=\"],\n  [\"18+39=\", \"34+31=\"],\n  [\"38+4=\", \"76-12=\"],\n  [\"64-29=\", \"30-21=\"],\n  [\"36-28=\", \"62-45=\"],\n  [\"82-3=\", \"91-55=\"],\n  [\"43-19=\", \"10+11=\"],\n  [\"59+20=\", \"43+37=\"],\n  [\"91+3=\", \"26+48=\"],\n  [\"35+5=\", \"96-71=\"],\n  [\"29+20=\", \"12+4=\"],\n  [\"71-21=\", \"37-34=\"],\n  [\"56+4=\", \"0+56=\"],\n  [\"99-31=\", \"70+15=\"],\n  [\"67-49=\", \"37-16=\"],\n  [\"37+6=\", \"39+21=\"],\n  [\"44+42=\", \"0+17=\"],\n  [\"40+29=\", \"9+38=\"],\n  [\"16-3=\", \"48+14=\"],\n  [\"95-7=\", \"31+47=\"],\n  [\"95-72=\", \"75-39=\"],\n  [\"76-33=\", \"65-2=\"],\n  [\"66-32=\", \"86-39=\"],\n  [\"84-13=\", \"14+66=\"],\n  [\"78-29=\", \"17+53=\"],\n  [\"19-6=\", \"36-22=\"],\n  [\"72-27=\", \"58+20=\"],\n  [\"64-12=\", \"36+42=\"],\n  [\"54+24=\", \"43+23=\"],\n  [\"50-41=\", \"56-4=\"],\n  [\"15-0=\", \"18+13=\"],\n  [\"5+83=\", \"32+10=\"],\n  [\"65-34=\", \"69-36=\"],\n  [\"1+43=\", \"32-14=\"],\n  [\"1+62=\", \"43+23=\"],\n  [\"90-52=\", \"35-6=\"],\n  [\"28-5=\", \"14+51=\"],\n  [\"15+53=\", \"41-12=\"],\n  [\"21+21=\", \"34-25=\"],\n  [\"44-30=\", \"99-49=\"],\n  [\"49+1=\", \"49-14=\"],\n  [\"65-28=\", \"54-11=\"],\n  [\"28+16=\", \"70-23=\"],\n  [\"87-44=\", \"45+37=\"],\n  [\"45+39=\", \"77+0=\"],\n  [\"98-28=\", \"50+2=\"],\n  [\"86-48=\", \"46-9=\"],\n  [\"71-34=\", \"27+16=\"],\n  [\"84-57=\", \"76+10=\"],\n  [\"1+52=\", \"2+85=\"],\n  [\"33+4=\", \"60-42=\"],\n  [\"1+57=\", \"58-57=\"],\n  [\"23+39=\", \"52-2=\"],\n  [\"86-58=\", \"65-25=\"],\n  [\"83-38=\", \"48-10=\"],\n  [\"97-8=\", \"48+44=\"],\n  [\"11+45=\", \"19+42=\"],\n  [\"81-58=\", \"8+20=\"],\n  [\"43+49=\", \"97-64=\"],\n  [\"71+6=\", \"63+19=\"],\n];\n\nconst body = context.document.body;\n\n// Phase 1: queue up a search for every old expression text.\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWildcards: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\n// Phase 2: replace each match (there is exactly one per expression) with\n// its new text.\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searchResults[i].items;\n  if (items.length > 0) {\n    items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace each arithmetic expression in the worksheet table with its\n# updated value. Every \"old\" expression text is unique in the document,\n# so Find/Replace on the whole document body is unambiguous and only\n# rewrites the matched run's text -- paragraph/run formatting (font,\n# size, alignment, etc.) is left untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2+39=', '66+2='),\n    @('81-14=', '36+0='),\n    @('61+29=', '42-0='),\n    @('26+31=', '94-35='),\n    @('15-7=', '48-24='),\n    @('62-5=', '15-1='),\n    @('70-14=', '19+37='),\n    @('62+12=', '43-1='),\n    @('84+10=', '44+9='),\n    @('31-18=', '45+34='),\n    @('25-23=', '93-74='),\n    @('6+83=', '26+62='),\n    @('88-69=', '68+14='),\n    @('53+11=', '24+7='),\n    @('69+12=', '46-15='),\n    @('7+86=', '3+8='),\n    @('59+4=', '38+41='),\n    @('40+39=', '32+20='),\n    @('76-28=', '82-76='),\n    @('27+27=', '93-26='),\n    @('5+78=', '49-11='),\n    @('10+55=', '73+26='),\n    @('60+29=', '81-40='),\n    @('73-69=', '41+20='),\n    @('41-30=', '72+8='),\n    @('40+53=', '36-24='),\n    @('5+3=', '63+16='),\n    @('67-48=', '1+55='),\n    @('4+81=', '92-86='),\n    @('93-21=', '74-4='),\n    @('77+6=', '6+3='),\n    @('31-23=', '2+63='),\n    @('18+26=', '72+24='),\n    @('70-54=', '86-51='),\n    @('18+60=', '16+43='),\n    @('82-62=', '80-63='),\n    @('56-46=', '76+8='),\n    @('90-66=', '57+8='),\n    @('58+30=', '74-25='),\n    @('67+17=', '85-59='),\n    @('71+22=', '3+86='),\n    @('18+39=', '34+31='),\n    @('38+4=', '76-12='),\n    @('64-29=', '30-21='),\n    @('36-28=', '62-45='),\n    @('82-3=', '91-55='),\n    @('43-19=', '10+11='),\n    @('59+20=', '43+37='),\n    @('91+3=', '26+48='),\n    @('35+5=', '96-71='),\n    @('29+20=', '12+4='),\n    @('71-21=', '37-34='),\n    @('56+4=', '0+56='),\n    @('99-31=', '70+15='),\n    @('67-49=', '37-16='),\n    @('37+6=', '39+21='),\n    @('44+42=', '0+17='),\n    @('40+29=', '9+38='),\n    @('16-3=', '48+14='),\n    @('95-7=', '31+47='),\n    @('95-72=', '75-39='),\n    @('76-33=', '65-2='),\n    @('66-32=', '86-39='),\n    @('84-13=', '14+66='),\n    @('78-29=', '17+53='),\n    @('19-6=', '36-22='),\n    @('72-27=', '58+20='),\n    @('64-12=', '36+42='),\n    @('54+24=', '43+23='),\n    @('50-41=', '56-4='),\n    @('15-0=', '18+13='),\n    @('5+83=', '32+10='),\n    @('65-34=', '69-36='),\n    @('1+43=', '32-14='),\n    @('1+62=', '43+23='),\n    @('90-52=', '35-6='),\n    @('28-5=', '14+51='),\n    @('15+53=', '41-12='),\n    @('21+21=', '34-25='),\n    @('44-30=', '99-49='),\n    @('49+1=', '49-14='),\n    @('65-28=', '54-11='),\n    @('28+16=', '70-23='),\n    @('87-44=', '45+37='),\n    @('45+39=', '77+0='),\n    @('98-28=', '50+2='),\n    @('86-48=', '46-9='),\n    @('71-34=', '27+16='),\n    @('84-57=', '76+10='),\n    @('1+52=', '2+85='),\n    @('33+4=', '60-42='),\n    @('1+57=', '58-57='),\n    @('23+39=', '52-2='),\n    @('86-58=', '65-25='),\n    @('83-38=', '48-10='),\n    @('97-8=', '48+44='),\n    @('11+45=', '19+42='),\n    @('81-58=', '8+20='),\n    @('43+49=', '97-64='),\n    @('71+6=', '63+19=')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
